{"js": "// The transcription markup for this page currently spells out the\n// opening \"<id>\" tag, the id value \"p158r_1\" and the closing \"</id>\"\n// tag as three separate runs. Collapse them into a single run whose\n// text is \"<id>p158r_1</id>\", keeping the formatting (Courier New /\n// dark-yellow / 9pt) that the first of those three runs already has.\nconst searchText = \"<id>p158r_1</id>\";\n\nconst results = context.document.body.search(searchText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the '<id>p158r_1</id>' run sequence to merge.\");\n}\n\n// Replacing the matched range with its own text causes Word to merge\n// the runs it spans into a single run, re-using the formatting of the\n// first run in the match (exactly mirroring what Word does when you\n// retype/replace text that already carries mixed run formatting).\nconst target = results.items[0];\ntarget.insertText(searchText, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# The transcription markup for this page currently spells out the\n# opening \"<id>\" tag, the id value \"p158r_1\" and the closing \"</id>\"\n# tag as three separate runs. Collapse them into a single run whose\n# text is \"<id>p158r_1</id>\", keeping the formatting (Courier New /\n# dark-yellow / 9pt) that the first of those three runs already has.\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Replacement.ClearFormatting()\n$rng.Find.Text = \"<id>p158r_1</id>\"\n$rng.Find.Replacement.Text = \"<id>p158r_1</id>\"\n$rng.Find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n"}
